# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.616.39'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '1.845.00'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').Value = "'0.9998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = "'314.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = "'0.9984"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').Value = "'0.4250"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.79%  '
$ws.Range('D8').Value = "'0.3643"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('D9').Value = "'0.07267"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('D10').Value = "'0.8948"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.09%  '
$ws.Range('D11').Value = "'20.63"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.08%  '
$ws.Range('D12').Value = '1.812.40'
$ws.Range('E12').Value = '  -3.31%  '
$ws.Range('D13').Value = "'6.584"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = "'5.362"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = "'0.06872"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = "'78.69"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').Value = "'0.000008864"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '27.602.50'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = "'4.986"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').Value = "'10.61"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').Value = '2.061.68'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').Value = "'2.030"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').Value = "'154.42"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'18.58"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = "'120.09"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.37%  '
$ws.Range('D29').Value = "'5.256"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').Value = "'1.839"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.64%  '
$ws.Range('D31').Value = "'0.08929"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('D32').Value = "'0.7816"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').Value = "'4.564"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.33%  '
$ws.Range('D34').Value = "'2.962"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').Value = "'1.106"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('D36').Value = "'0.9979"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = "'1.102"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.05403"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('D39').Value = "'0.01935"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = "'2.800"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('D41').Value = "'6.890"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').Value = "'0.5069"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.19%  '
$ws.Range('D43').Value = "'0.1650"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('D44').Value = "'8.256"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.24%  '
$ws.Range('D45').Value = "'0.06628"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('D46').Value = "'10.36"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('D47').Value = "'0.4715"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('D48').Value = "'105.11"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = "'0.9978"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').Value = "'1.805"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.47%  '
